# Generate Report for handback
# Row 3 of the "zh-cn" and "de-de" status sheets previously carried
# placeholder handoff/handback timestamps copied from row 2. This fills
# in the real, distinct handoff/handback datetimes for the
# e716d8a8-7d6e-492c-b3ce-9abf02466350 entry on each language sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-13 11:36:06"
$wsZhCn.Range("G3").Value = "2016-01-13 11:37:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-13 11:36:30"
$wsDeDe.Range("G3").Value = "2016-01-13 11:38:15"
